# Update "Datos actualizados ..." timestamp cell (A1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 01:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 737865
$ws.Range("C4").Value = 28130
$ws.Range("D4").Value = 67483
$ws.Range("E4").Value = 631422
$ws.Range("F4").Value = 13551
$ws.Range("G4").Value = 1806
$ws.Range("H4").Value = 38960

# Row 5 - Espana
$ws.Range("B5").Value = 194416
$ws.Range("C5").Value = 3577
$ws.Range("E5").Value = 98980

# Row 8 - Alemania
$ws.Range("B8").Value = 143475
$ws.Range("C8").Value = 2078
$ws.Range("E8").Value = 53598
$ws.Range("G8").Value = 125
$ws.Range("H8").Value = 4477

# Row 15 - Brasil
$ws.Range("B15").Value = 36722
$ws.Range("C15").Value = 3040
$ws.Range("E15").Value = 20335
$ws.Range("G15").Value = 220
$ws.Range("H15").Value = 2361

# Row 55 - Argentina
$ws.Range("B55").Value = 2839
$ws.Range("C55").Value = 170
$ws.Range("E55").Value = 2022
$ws.Range("G55").Value = 9
$ws.Range("H55").Value = 132

# Row 100 - Uruguay
$ws.Range("B100").Value = 517
$ws.Range("C100").Value = 15
$ws.Range("D100").Value = 298
$ws.Range("E100").Value = 210
$ws.Range("F100").Value = 14

# Row 121 - Venezuela
$ws.Range("D121").Value = 117
$ws.Range("E121").Value = 101

# Row 172 - Nepal
$ws.Range("D172").Value = 3
$ws.Range("E172").Value = 28

# Row 205 - Mauritania
$ws.Range("D205").Value = 6
$ws.Range("E205").Value = 0
